{"js": "// Apply the five text replacements described by the diff.\n// Each target string is unique in the document, so a plain text search\n// (no wildcards) finds exactly one range per call.\nconst replacements = [\n  [\"CV \\u2013 Patti Fernandez\", \"CV : Patti Fernandez\"],\n  [\n    \"ABC Studios\\u00a0: Lead Animator (Jan 2018 - Present)\",\n    \"ABC Studios : Animateur principal (depuis janvier 2018)\",\n  ],\n  [\n    \"XYZ Media\\u00a0: Senior Animator (juin 2015 - Dec 2017)\",\n    \"XYZ Media : animateur senior (juin 2015 \\u00e0 d\\u00e9cembre 2017)\",\n  ],\n  [\n    \"MNO Entertainment\\u00a0: Junior Animator (sep 2012 - mai 2015)\",\n    \"MNO Entertainment : Animateur junior (septembre 2012 \\u00e0 mai 2015)\",\n  ],\n  [\n    \"L\\u2019art de l\\u2019animation : Guide pour les d\\u00e9butants.\",\n    \"The Art of 3D Animation: A Guide for Beginners.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the five text replacements described by the diff using Find/Replace\n# on the document's main Range. Each target string occurs exactly once, so\n# wdReplaceOne keeps the replacement scoped to that single occurrence.\n\n$d = $word.ActiveDocument\n\nfunction Replace-OnceInDoc($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\nReplace-OnceInDoc \"CV \u2013 Patti Fernandez\" \"CV : Patti Fernandez\"\nReplace-OnceInDoc \"ABC Studios : Lead Animator (Jan 2018 - Present)\" \"ABC Studios : Animateur principal (depuis janvier 2018)\"\nReplace-OnceInDoc \"XYZ Media : Senior Animator (juin 2015 - Dec 2017)\" \"XYZ Media : animateur senior (juin 2015 \u00e0 d\u00e9cembre 2017)\"\nReplace-OnceInDoc \"MNO Entertainment : Junior Animator (sep 2012 - mai 2015)\" \"MNO Entertainment : Animateur junior (septembre 2012 \u00e0 mai 2015)\"\nReplace-OnceInDoc \"L\u2019art de l\u2019animation : Guide pour les d\u00e9butants.\" \"The Art of 3D Animation: A Guide for Beginners.\"\n"}
